$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '28.867.81'
$ws.Cells.Item(2, 5).Value = '  -1.43%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.904.41'
$ws.Cells.Item(3, 5).Value = '  -1.47%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.18%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '324.58'
$ws.Cells.Item(5, 5).Value = '  -0.25%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.83%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.3813'

# Row 9
$ws.Cells.Item(9, 4).Value = '0.07714'
$ws.Cells.Item(9, 5).Value = '  -1.32%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.9785'

# Row 11
$ws.Cells.Item(11, 4).Value = '22.19'
$ws.Cells.Item(11, 5).Value = '  -2.02%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.880.16'
$ws.Cells.Item(12, 5).Value = '  -2.20%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Chainlink'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(13, 4).Value = '6.946'
$ws.Cells.Item(13, 5).Value = '  -1.91%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = '5.671'
$ws.Cells.Item(14, 5).Value = '  -2.13%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '0.07065'
$ws.Cells.Item(15, 5).Value = '  +0.08%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.08%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '83.81'
$ws.Cells.Item(17, 5).Value = '  -3.44%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000009450'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -2.95%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '16.62'
$ws.Cells.Item(19, 5).Value = '  -2.34%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.01%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '28.857.04'
$ws.Cells.Item(21, 5).Value = '  -1.44%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -2.98%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -1.17%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.100'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.34%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '158.82'
$ws.Cells.Item(25, 5).Value = '  +0.80%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '18.99'
$ws.Cells.Item(26, 5).Value = '  -1.88%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '5.672'
$ws.Cells.Item(27, 5).Value = '  -1.76%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '117.53'
$ws.Cells.Item(28, 5).Value = '  -1.25%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '1.872'
$ws.Cells.Item(29, 5).Value = '  +1.58%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.09300'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.30%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '0.8614'
$ws.Cells.Item(31, 5).Value = '  -0.35%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '5.094'
$ws.Cells.Item(32, 5).Value = '  -1.33%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.243'
$ws.Cells.Item(33, 5).Value = '  -4.33%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '3.029'
$ws.Cells.Item(34, 5).Value = '  -1.68%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.05706'
$ws.Cells.Item(35, 5).Value = '  -1.24%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.160'
$ws.Cells.Item(36, 4).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.13%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.02041'
$ws.Cells.Item(38, 5).Value = '  -1.83%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '7.455'
$ws.Cells.Item(39, 5).Value = '  -2.43%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '0.5481'
$ws.Cells.Item(40, 5).Value = '  -3.11%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'MXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(41, 4).Value = '2.911'
$ws.Cells.Item(41, 5).Value = '  +7.39%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).Value = '0.1751'
$ws.Cells.Item(42, 5).Value = '  -1.71%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '9.329'

# Row 44
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 4).Value = '2.169'
$ws.Cells.Item(44, 5).Value = '  +4.67%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'PEPE'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.000002738'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -10.70%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.5174'
$ws.Cells.Item(46, 5).Value = '  -1.54%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.73%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '0.06881'
$ws.Cells.Item(48, 5).Value = '  +0.33%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '1.777'
$ws.Cells.Item(49, 5).Value = '  -2.00%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '110.44'
$ws.Cells.Item(50, 5).Value = '  -0.66%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.2859'
$ws.Cells.Item(51, 5).Value = '  -4.23%  '
